$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 182.44444
$ws.Range("I33").Value = 114
$ws.Range("K33").Value = 114
$ws.Range("M33").Value = 115

$ws.Range("H41").Value = 548.5833
$ws.Range("I41").Value = 775
$ws.Range("J41").Value = 435.375
$ws.Range("K41").Value = 775
$ws.Range("L41").Value = 435.375
$ws.Range("M41").Value = -335
$ws.Range("N41").Value = -1315.375

$ws.Range("H51").Value = 5297
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 5297
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 5297
$ws.Range("N51").Value = -6265
$ws.Range("M51").ClearContents()

$ws.Range("H111").Value = 415.1111
$ws.Range("I111").Value = 157.25
$ws.Range("J111").Value = 621.4
$ws.Range("K111").Value = 471.75
$ws.Range("L111").Value = 1864.2
$ws.Range("M111").Value = 2595.25
$ws.Range("N111").Value = -7998.2

$ws.Range("H116").Value = 5177.067
$ws.Range("I116").Value = 1328.6666
$ws.Range("J116").Value = 10949.667
$ws.Range("K116").Value = 1328.6666
$ws.Range("L116").Value = 10949.667
$ws.Range("M116").Value = 2113.3334
$ws.Range("N116").Value = -17833.667

$ws.Range("H129").Value = 833.8431399999999
$ws.Range("J129").Value = 955.475
$ws.Range("L129").Value = 2866.425
$ws.Range("N129").Value = -12866.425

$ws.Range("H137").Value = 6384750.5
$ws.Range("I137").Value = 20001038
$ws.Range("J137").Value = 2115.8438
$ws.Range("K137").Value = 60003114
$ws.Range("L137").Value = 6347.5314
$ws.Range("M137").Value = -60000564
$ws.Range("N137").Value = -11447.5314

$ws.Range("H138").Value = 3625121.5
$ws.Range("I138").Value = 1275.762
$ws.Range("J138").Value = 6669152
$ws.Range("K138").Value = 3827.286
$ws.Range("L138").Value = 20007456
$ws.Range("M138").Value = 1312.714
$ws.Range("N138").Value = -20017736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6670084
$ws.Range("I32").Value = 9262505
$ws.Range("J32").Value = 3859.2856
$ws.Range("K32").Value = 9262505
$ws.Range("L32").Value = 3859.2856
$ws.Range("M32").Value = -9262218
$ws.Range("N32").Value = -4433.2856

$ws.Range("H61").Value = 250500620
$ws.Range("I61").Value = 500500000
$ws.Range("J61").Value = 501250
$ws.Range("K61").Value = 500500000
$ws.Range("L61").Value = 501250
$ws.Range("M61").Value = -500499788
$ws.Range("N61").Value = -501674

$ws.Range("H74").Value = 13264830
$ws.Range("I74").Value = 22819860
$ws.Range("J74").Value = 126664.5
$ws.Range("K74").Value = 22819860
$ws.Range("L74").Value = 126664.5
$ws.Range("M74").Value = -22818986
$ws.Range("N74").Value = -128412.5

$ws.Range("H77").Value = 13264830
$ws.Range("I77").Value = 22819860
$ws.Range("J77").Value = 126664.5
$ws.Range("K77").Value = 114099300
$ws.Range("L77").Value = 633322.5
$ws.Range("M77").Value = -114094932
$ws.Range("N77").Value = -642058.5

$ws.Range("H136").Value = 250500620
$ws.Range("I136").Value = 500500000
$ws.Range("J136").Value = 501250
$ws.Range("K136").Value = 1501500000
$ws.Range("L136").Value = 1503750
$ws.Range("M136").Value = -1501497450
$ws.Range("N136").Value = -1508850

$ws.Range("H137").Value = 49933.332
$ws.Range("J137").Value = 49933.332
$ws.Range("L137").Value = 49933.332
$ws.Range("N137").Value = -60133.332

$ws.Range("H139").Value = 59715
$ws.Range("J139").Value = 59715
$ws.Range("L139").Value = 59715
$ws.Range("N139").Value = -69995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1054.1111
$ws.Range("I20").Value = 1139.6
$ws.Range("J20").Value = 947.25
$ws.Range("K20").Value = 1139.6
$ws.Range("L20").Value = 947.25
$ws.Range("M20").Value = -892.5999999999999
$ws.Range("N20").Value = -1441.25

$ws.Range("H134").Value = 8279.429
$ws.Range("I134").Value = 7568
$ws.Range("K134").Value = 22704
$ws.Range("M134").Value = -20169

$ws.Range("H138").Value = 32800
$ws.Range("J138").Value = 32800
$ws.Range("L138").Value = 32800
$ws.Range("N138").Value = -43080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 20818670
$ws.Range("I58").Value = 31877946
$ws.Range("J58").Value = 1211.5294
$ws.Range("K58").Value = 31877946
$ws.Range("L58").Value = 1211.5294
$ws.Range("M58").Value = -31877743
$ws.Range("N58").Value = -1617.5294

$ws.Range("H132").Value = 41252.152
$ws.Range("J132").Value = 93628.37
$ws.Range("L132").Value = 280885.11
$ws.Range("N132").Value = -285945.11

$ws.Range("H136").Value = 20818670
$ws.Range("I136").Value = 31877946
$ws.Range("J136").Value = 1211.5294
$ws.Range("K136").Value = 95633838
$ws.Range("L136").Value = 3634.5882
$ws.Range("M136").Value = -95631288
$ws.Range("N136").Value = -8734.5882

$ws.Range("H138").Value = 49800
$ws.Range("J138").Value = 49800
$ws.Range("L138").Value = 49800
$ws.Range("N138").Value = -60080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11819235
$ws.Range("I4").Value = 5000598
$ws.Range("J4").Value = 20001600
$ws.Range("K4").Value = 15001794
$ws.Range("L4").Value = 60004800
$ws.Range("M4").Value = -15001682
$ws.Range("N4").Value = -60005024

$ws.Range("H17").Value = 1033
$ws.Range("I17").Value = 1043.3334
$ws.Range("J17").Value = 1002
$ws.Range("K17").Value = 3130.0002
$ws.Range("L17").Value = 3006
$ws.Range("M17").Value = -2961.0002
$ws.Range("N17").Value = -3344

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 76842.86
$ws.Range("I70").Value = 105080
$ws.Range("J70").Value = 6250
$ws.Range("K70").Value = 105080
$ws.Range("L70").Value = 6250
$ws.Range("M70").Value = -104810
$ws.Range("N70").Value = -6790

$ws.Range("H73").Value = 76842.86
$ws.Range("I73").Value = 105080
$ws.Range("J73").Value = 6250
$ws.Range("K73").Value = 105080
$ws.Range("L73").Value = 6250
$ws.Range("M73").Value = -104144
$ws.Range("N73").Value = -8122

$ws.Range("H124").Value = 50780
$ws.Range("J124").Value = 50780
$ws.Range("L124").Value = 50780
$ws.Range("N124").Value = -60600

$ws.Range("H132").Value = 61740.707
$ws.Range("I132").Value = 65520.5
$ws.Range("J132").Value = 58380.89
$ws.Range("K132").Value = 196561.5
$ws.Range("L132").Value = 175142.67
$ws.Range("M132").Value = -194031.5
$ws.Range("N132").Value = -180202.67

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 259.66666
$ws.Range("I55").Value = 237.64285
$ws.Range("J55").Value = 290.5
$ws.Range("K55").Value = 237.64285
$ws.Range("L55").Value = 290.5
$ws.Range("M55").Value = -64.64285000000001
$ws.Range("N55").Value = -636.5

$ws.Range("H136").Value = 85168.625
$ws.Range("I136").Value = 45304.22
$ws.Range("J136").Value = 161575.42
$ws.Range("K136").Value = 135912.66
$ws.Range("L136").Value = 484726.26
$ws.Range("M136").Value = -133362.66
$ws.Range("N136").Value = -489826.26

$ws.Range("H139").Value = 51215
$ws.Range("J139").Value = 51215
$ws.Range("L139").Value = 51215
$ws.Range("N139").Value = -61495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 45000
$ws.Range("J87").Value = 45000
$ws.Range("L87").Value = 45000
$ws.Range("N87").Value = -47496

$ws.Range("H90").Value = 45000
$ws.Range("J90").Value = 45000
$ws.Range("L90").Value = 135000
$ws.Range("N90").Value = -147480

$ws.Range("H136").Value = 54986.242
$ws.Range("I136").Value = 39351.19
$ws.Range("J136").Value = 91941.82000000001
$ws.Range("K136").Value = 118053.57
$ws.Range("L136").Value = 275825.46
$ws.Range("M136").Value = -115503.57
$ws.Range("N136").Value = -280925.46
